$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 717.431
$ws.Range("J17").Value = 723.7544
$ws.Range("L17").Value = 2171.2632
$ws.Range("N17").Value = -2507.2632
$ws.Range("H33").Value = 341.5263
$ws.Range("I33").Value = 158.8125
$ws.Range("J33").Value = 1316
$ws.Range("K33").Value = 158.8125
$ws.Range("L33").Value = 1316
$ws.Range("M33").Value = 70.1875
$ws.Range("N33").Value = -1774
$ws.Range("H53").Value = 422.76923
$ws.Range("I53").Value = 81.7
$ws.Range("K53").Value = 81.7
$ws.Range("M53").Value = 555.3
$ws.Range("H58").Value = 1942.4286
$ws.Range("I58").Value = 1900
$ws.Range("J58").Value = 1999
$ws.Range("K58").Value = 5700
$ws.Range("L58").Value = 5997
$ws.Range("M58").Value = -5550
$ws.Range("N58").Value = -6297
$ws.Range("H76").Value = 5343.4
$ws.Range("I76").Value = 4491.75
$ws.Range("J76").Value = 8750
$ws.Range("K76").Value = 4491.75
$ws.Range("L76").Value = 8750
$ws.Range("M76").Value = -4176.75
$ws.Range("N76").Value = -9380
$ws.Range("H79").Value = 5343.4
$ws.Range("I79").Value = 4491.75
$ws.Range("J79").Value = 8750
$ws.Range("K79").Value = 4491.75
$ws.Range("L79").Value = 8750
$ws.Range("M79").Value = -3399.75
$ws.Range("N79").Value = -10934
$ws.Range("H86").Value = 4745.2666
$ws.Range("I86").Value = 4802.727
$ws.Range("J86").Value = 4587.25
$ws.Range("K86").Value = 4802.727
$ws.Range("L86").Value = 4587.25
$ws.Range("M86").Value = -3679.727
$ws.Range("N86").Value = -6833.25
$ws.Range("H89").Value = 4745.2666
$ws.Range("I89").Value = 4802.727
$ws.Range("J89").Value = 4587.25
$ws.Range("K89").Value = 24013.635
$ws.Range("L89").Value = 22936.25
$ws.Range("M89").Value = -18397.635
$ws.Range("N89").Value = -34168.25
$ws.Range("H98").Value = 2113.742
$ws.Range("I98").Value = 2158.25
$ws.Range("K98").Value = 2158.25
$ws.Range("M98").Value = -660.25
$ws.Range("H111").Value = 1198
$ws.Range("I111").Value = 1024.5
$ws.Range("K111").Value = 3073.5
$ws.Range("M111").Value = -6.5
$ws.Range("H122").Value = 2113.742
$ws.Range("I122").Value = 2158.25
$ws.Range("K122").Value = 6474.75
$ws.Range("M122").Value = -4024.75
$ws.Range("H132").Value = 9393.542
$ws.Range("I132").Value = 9185.25
$ws.Range("K132").Value = 27555.75
$ws.Range("M132").Value = -25025.75
$ws.Range("H135").Value = 6008
$ws.Range("I135").Value = 11793.75
$ws.Range("J135").Value = 1379.4
$ws.Range("K135").Value = 106143.75
$ws.Range("L135").Value = 12414.6
$ws.Range("M135").Value = -103608.75
$ws.Range("N135").Value = -17484.6
$ws.Range("H137").Value = 3217.261
$ws.Range("I137").Value = 2596.0833
$ws.Range("K137").Value = 7788.249899999999
$ws.Range("M137").Value = -5238.249899999999
$ws.Range("H138").Value = 21928.541
$ws.Range("I138").Value = 25727.115
$ws.Range("K138").Value = 77181.345
$ws.Range("M138").Value = -72041.345
$ws.Range("H141").Value = 1375.5
$ws.Range("I141").Value = 1333.8667
$ws.Range("K141").Value = 4001.6001
$ws.Range("M141").Value = 1178.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14047.671
$ws.Range("I32").Value = 14090.561
$ws.Range("K32").Value = 14090.561
$ws.Range("M32").Value = -13803.561
$ws.Range("H74").Value = 1306.8572
$ws.Range("I74").Value = 1292
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1292
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -418
$ws.Range("N74").Value = -3248
$ws.Range("H77").Value = 1306.8572
$ws.Range("I77").Value = 1292
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 6460
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -2092
$ws.Range("N77").Value = -16236
$ws.Range("H97").Value = 3139.9
$ws.Range("J97").Value = 3800
$ws.Range("L97").Value = 3800
$ws.Range("N97").Value = -4792
$ws.Range("H110").Value = 2107.2222
$ws.Range("I110").Value = 1119.5
$ws.Range("K110").Value = 1119.5
$ws.Range("M110").Value = 925.5
$ws.Range("H132").Value = 75652.93
$ws.Range("I132").Value = 114688.11
$ws.Range("K132").Value = 344064.33
$ws.Range("M132").Value = -341534.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 8957.8
$ws.Range("I11").Value = 445
$ws.Range("J11").Value = 14633
$ws.Range("K11").Value = 445
$ws.Range("L11").Value = 14633
$ws.Range("M11").Value = -305
$ws.Range("N11").Value = -14913
$ws.Range("H128").Value = 1800
$ws.Range("I128").Value = 1800
$ws.Range("K128").Value = 5400
$ws.Range("M128").Value = -2910

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 100.666664
$ws.Range("I2").Value = 80.8
$ws.Range("K2").Value = 80.8
$ws.Range("M2").Value = 32.2
$ws.Range("H86").Value = 6249.25
$ws.Range("I86").Value = 6249.25
$ws.Range("K86").Value = 6249.25
$ws.Range("M86").Value = -5126.25
$ws.Range("H89").Value = 6249.25
$ws.Range("I89").Value = 6249.25
$ws.Range("K89").Value = 31246.25
$ws.Range("M89").Value = -25630.25
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = $null
$ws.Range("N94").Value = $null
$ws.Range("H134").Value = 169166.5
$ws.Range("I134").Value = 501999.5
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 1505998.5
$ws.Range("L134").Value = 8250
$ws.Range("M134").Value = -1503463.5
$ws.Range("N134").Value = -13320
$ws.Range("H138").Value = 79998.164
$ws.Range("J138").Value = 79998.164
$ws.Range("L138").Value = 79998.164
$ws.Range("N138").Value = -90278.164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 491.85715
$ws.Range("I121").Value = 163.33333
$ws.Range("K121").Value = 489.99999
$ws.Range("M121").Value = 820.00001
$ws.Range("H129").Value = 503094.12
$ws.Range("J129").Value = 853464.3
$ws.Range("L129").Value = 2560392.9
$ws.Range("N129").Value = -2570392.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 501999
$ws.Range("I132").Value = 501999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1505997
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1503467
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10354.728
$ws.Range("I40").Value = 5975.75
$ws.Range("J40").Value = 12857
$ws.Range("K40").Value = 5975.75
$ws.Range("L40").Value = 12857
$ws.Range("M40").Value = -5839.75
$ws.Range("N40").Value = -13129
$ws.Range("H55").Value = 1122.7
$ws.Range("I55").Value = 576
$ws.Range("J55").Value = 1669.4
$ws.Range("K55").Value = 576
$ws.Range("L55").Value = 1669.4
$ws.Range("M55").Value = -403
$ws.Range("N55").Value = -2015.4
$ws.Range("H68").Value = 8278.4
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 8278.4
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 8278.4
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = -9776.4
$ws.Range("H71").Value = 8278.4
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 8278.4
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 41392
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = -48880
$ws.Range("H75").Value = 50160
$ws.Range("I75").Value = 50160
$ws.Range("K75").Value = 50160
$ws.Range("M75").Value = -49224
$ws.Range("H78").Value = 50160
$ws.Range("I78").Value = 50160
$ws.Range("K78").Value = 150480
$ws.Range("M78").Value = -145800
$ws.Range("H82").Value = 2618.7334
$ws.Range("I82").Value = 2498.2
$ws.Range("K82").Value = 2498.2
$ws.Range("M82").Value = -2137.2
$ws.Range("H85").Value = 2618.7334
$ws.Range("I85").Value = 2498.2
$ws.Range("K85").Value = 2498.2
$ws.Range("M85").Value = -1250.2
$ws.Range("H93").Value = 2852.889
$ws.Range("I93").Value = 2852.889
$ws.Range("K93").Value = 2852.889
$ws.Range("M93").Value = -1604.889
$ws.Range("H100").Value = 2060.4546
$ws.Range("I100").Value = 1096.1875
$ws.Range("J100").Value = 4631.8335
$ws.Range("K100").Value = 1096.1875
$ws.Range("L100").Value = 4631.8335
$ws.Range("M100").Value = -555.1875
$ws.Range("N100").Value = -5713.8335
$ws.Range("H132").Value = 53486.418
$ws.Range("I132").Value = 63008.7
$ws.Range("J132").Value = 5875
$ws.Range("K132").Value = 189026.1
$ws.Range("L132").Value = 17625
$ws.Range("M132").Value = -186496.1
$ws.Range("N132").Value = -22685
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 340499.66
$ws.Range("I62").Value = 10500
$ws.Range("K62").Value = 10500
$ws.Range("M62").Value = -9876
$ws.Range("H65").Value = 340499.66
$ws.Range("I65").Value = 10500
$ws.Range("K65").Value = 52500
$ws.Range("M65").Value = -49380
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = $null
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null
$ws.Range("H124").Value = 73919
$ws.Range("J124").Value = 73919
$ws.Range("L124").Value = 73919
$ws.Range("N124").Value = -83739
$ws.Range("H126").Value = 63521.766
$ws.Range("I126").Value = 75632.86
$ws.Range("J126").Value = 7003.3335
$ws.Range("K126").Value = 226898.58
$ws.Range("L126").Value = 21010.0005
$ws.Range("M126").Value = -224428.58
$ws.Range("N126").Value = -25950.0005
